$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.627.89'
$ws.Range('E2').Value = '  -1.80%  '
$ws.Range('D3').Value = '3.072.32'
$ws.Range('E3').Value = '  -2.29%  '
$ws.Range('E4').Value = '  -0.69%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '593.32'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '155.13'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +1.44%  '
$ws.Range('E7').Value = '  -0.28%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.536'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.40%  '
$ws.Range('D9').Value = '3.071.08'
$ws.Range('E9').Value = '  -2.23%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.157'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -2.69%  '
$ws.Range('E11').Value = '  -1.82%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.451'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -2.98%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000237'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -3.94%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.55'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -4.35%  '
$ws.Range('E15').Value = '  +0.27%  '
$ws.Range('D16').Value = '3.578.17'
$ws.Range('E16').Value = '  -2.44%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.18'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.14%  '
$ws.Range('D18').Value = '63.496.78'
$ws.Range('E18').Value = '  -1.34%  '
$ws.Range('D19').Value = '3.069.93'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '480.25'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +1.31%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.39'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -4.20%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.708'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -5.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.53'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -1.78%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.41'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.89%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '81.30'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.16%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.80'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -4.43%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.60'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +5.82%  '
$ws.Range('E28').Value = '  +0.07%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.57'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +1.93%  '
$ws.Range('E30').Value = '  -1.66%  '
$ws.Range('E31').Value = '  -0.74%  '
$ws.Range('E32').Value = '  -2.28%  '
$ws.Range('E33').Value = '  -5.23%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '27.14'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.86%  '
$ws.Range('D35').Value = '0.0₃0834'
$ws.Range('E35').Value = '  -4.57%  '
$ws.Range('E36').Value = '  -0.55%  '
$ws.Range('B37').Value = 'dogwifhat'
$ws.Range('C37').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.31'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -3.62%  '
$ws.Range('B38').Value = 'Filecoin'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.99'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -3.59%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.23'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -3.59%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '50.76'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.91%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '9.22'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -1.89%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '437.52'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -5.83%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.290'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.92%  '
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0362'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -4.62%  '
$ws.Range('B45').Value = 'Kaspa'
$ws.Range('C45').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.112'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.75%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '40.02'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.93%  '
$ws.Range('D47').Value = '2.818.53'
$ws.Range('E47').Value = '  -2.56%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '131.63'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.49%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '25.39'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -1.53%  '
$ws.Range('E50').Value = '  +0.01%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.24'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -3.06%  '
